# Apply updated values after correcting error estimation and number of
# projected years (trends not yet rerun after filtering bug fix).

$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("C2").Value = 5
$ws1.Range("E2").Value = 20.8
$ws1.Range("C3").Value = 10
$ws1.Range("E3").Value = 41.7
$ws1.Range("E4").Value = 33.3
$ws1.Range("C5").Value = 0
$ws1.Range("E5").Value = 0
$ws1.Range("C6").Value = 1
$ws1.Range("E6").Value = 4.2
$ws1.Range("C7").Value = 25

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("C4").Value = 24

# --- Sheet: "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")
$ws5.Range("B2").Value = 87
$ws5.Range("C2").Value = 84.5
$ws5.Range("D2").Value = 87
$ws5.Range("E2").Value = 98.90000000000001
$ws5.Range("B4").Value = 13
$ws5.Range("C4").Value = 12.6
$ws5.Range("D4").Value = 1
$ws5.Range("E4").Value = 1.1

$wb.Save()
